$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 1293 -> 1294, F6 62 -> 63
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1294
$wsExhibit.Range("F6").Value = 63

# Sheet "全部类型": F5 1293 -> 1294, F7 62 -> 63
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1294
$wsAll.Range("F7").Value = 63
